# Update "想去人数" (want-to-go count) figures pulled from the latest
# bilibili scrape on sheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3409
$wsExhibit.Range("F4").Value = 134
$wsExhibit.Range("F5").Value = 6967
$wsExhibit.Range("F6").Value = 2431
$wsExhibit.Range("F7").Value = 41
$wsExhibit.Range("F8").Value = 109
$wsExhibit.Range("F14").Value = 570

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3409
$wsAll.Range("F5").Value = 134
$wsAll.Range("F6").Value = 6967
$wsAll.Range("F7").Value = 2431
$wsAll.Range("F8").Value = 41
$wsAll.Range("F9").Value = 109
$wsAll.Range("F15").Value = 570
